# Add an "ema" (7-period exponential moving average of Close) column in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (bold/centered/bordered style) from F1 onto G1,
# then set its label.
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Cells.Item(1, 7).Value2 = "ema"

# EMA with period 7 -> smoothing factor alpha = 2 / (period + 1) = 0.25.
# Seed the series with the first Close value (row 2), as is conventional.
$alpha = 0.25

$prevEma = $ws.Cells.Item(2, 5).Value2
$ws.Cells.Item(2, 7).Value2 = $prevEma

for ($r = 3; $r -le 365; $r++) {
    $close = $ws.Cells.Item($r, 5).Value2
    $ema = $alpha * $close + (1 - $alpha) * $prevEma
    $ws.Cells.Item($r, 7).Value2 = $ema
    $prevEma = $ema
}
